$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $ws.Range($cell).Style
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $value
    $ws.Range($cell).Style = $origStyle
}

Set-TextValue 'D2' '56.623.02'
Set-TextValue 'E2' '  -0.24%  '
Set-TextValue 'D3' '2.319.17'
Set-TextValue 'E3' '  -0.29%  '
Set-TextValue 'E4' '  -0.09%  '
Set-TextValue 'D5' '515.20'
Set-TextValue 'E5' '  -1.18%  '
Set-TextValue 'D6' '131.90'
Set-TextValue 'E6' '  -2.51%  '
Set-TextValue 'E7' '  +0.39%  '
Set-TextValue 'D8' '0.533'
Set-TextValue 'E8' '  -0.76%  '
Set-TextValue 'D9' '0.100'
Set-TextValue 'E9' '  -3.08%  '
Set-TextValue 'E10' '  -0.09%  '
Set-TextValue 'E11' '  -1.47%  '
Set-TextValue 'E12' '  -2.07%  '
Set-TextValue 'D13' '23.61'
Set-TextValue 'E13' '  -1.93%  '
Set-TextValue 'D14' '2.732.60'
Set-TextValue 'E14' '  -0.31%  '
Set-TextValue 'D15' '56.606.26'
Set-TextValue 'E15' '  -0.37%  '
Set-TextValue 'E16' '  -1.52%  '
Set-TextValue 'D17' '2.318.78'
Set-TextValue 'E17' '  -0.83%  '
Set-TextValue 'E18' '  -1.82%  '
Set-TextValue 'D19' '328.21'
Set-TextValue 'E19' '  +1.39%  '
Set-TextValue 'D20' '4.14'
Set-TextValue 'E20' '  -2.24%  '
Set-TextValue 'E21' '  +1.59%  '
Set-TextValue 'E22' '  -0.16%  '
Set-TextValue 'D23' '61.08'
Set-TextValue 'E23' '  +0.68%  '
Set-TextValue 'E24' '  -0.88%  '
Set-TextValue 'D25' '8.60'
Set-TextValue 'E25' '  +7.28%  '
Set-TextValue 'E26' '  +0.55%  '
Set-TextValue 'E27' '  +0.65%  '
Set-TextValue 'D28' '167.65'
Set-TextValue 'E28' '  +0.40%  '
Set-TextValue 'E29' '  -2.53%  '
Set-TextValue 'D30' '0.0₃0717'
Set-TextValue 'E30' '  -4.05%  '
Set-TextValue 'E31' '  -2.28%  '
Set-TextValue 'D32' '18.28'
Set-TextValue 'E32' '  -0.67%  '
Set-TextValue 'E33' '  -0.01%  '
Set-TextValue 'D34' '0.997'
Set-TextValue 'E34' '  +0.52%  '
Set-TextValue 'E35' '  -1.52%  '
Set-TextValue 'D36' '3.93'
Set-TextValue 'E36' '  -3.11%  '
Set-TextValue 'D37' '0.881'
Set-TextValue 'E37' '  -5.40%  '
Set-TextValue 'E38' '  +0.29%  '
Set-TextValue 'D39' '38.63'
Set-TextValue 'E39' '  +1.79%  '
Set-TextValue 'D40' '148.54'
Set-TextValue 'E40' '  +6.86%  '
Set-TextValue 'E41' '  -1.52%  '
Set-TextValue 'D42' '3.55'
Set-TextValue 'E42' '  -1.83%  '
Set-TextValue 'D43' '275.43'
Set-TextValue 'E43' '  -1.66%  '
Set-TextValue 'E44' '  -5.80%  '
Set-TextValue 'D45' '0.0929'
Set-TextValue 'E45' '  -0.58%  '
Set-TextValue 'D46' '0.0494'
Set-TextValue 'E46' '  -2.57%  '
Set-TextValue 'E47' '  -1.48%  '
Set-TextValue 'D48' '18.18'
Set-TextValue 'E48' '  -0.02%  '
Set-TextValue 'B49' 'Polygon'
Set-TextValue 'C49' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D49' '0.380'
Set-TextValue 'E49' '  +0.10%  '
Set-TextValue 'B50' 'VeChain'
Set-TextValue 'C50' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D50' '0.0214'
Set-TextValue 'E50' '  -2.05%  '
Set-TextValue 'D51' '17.05'
Set-TextValue 'E51' '  +0.51%  '
